# Hasil dan Pembahasan(AutoRecovered).xlsx - "Add files via upload" edit
# Fills in the Recall / Precision / F-Measure rows on the "Evaluasi" sheet
# with live formulas (formatted as percentages) and updates the saved
# selection, mirroring the changes captured in the authoritative diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluasi")

# --- Row 7: Recall = TP / (row total) ------------------------------------
$ws.Range("F7").Formula = "=F4/5"
$ws.Range("G7").Formula = "=G5/5"
$ws.Range("H7").Formula = "=H6/5"
$ws.Range("I7").Formula = "=I4/5"
$ws.Range("J7").Formula = "=J5/5"
$ws.Range("K7").Formula = "=K6/5"
$ws.Range("L7").Formula = "=L4/5"
$ws.Range("M7").Formula = "=M5/5"
$ws.Range("N7").Formula = "=N6/5"

# --- Row 8: Precision = TP / (column total) -------------------------------
$ws.Range("F8").Formula = "=F4/4"
$ws.Range("G8").Formula = "=5/6"
$ws.Range("H8").Formula = "=5/5"
$ws.Range("I8").Formula = "=I4/4"
$ws.Range("J8").Formula = "=5/6"
$ws.Range("K8").Formula = "=5/5"
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 1

# --- Row 9: F-Measure = 2 * (Recall * Precision) / (Recall + Precision) --
$ws.Range("F9").Formula = "=2*((F7*F8)/(F7+F8))"
$ws.Range("G9").Formula = "=2*((G7*G8)/(G7+G8))"
$ws.Range("H9").Formula = "=2*((H7*H8)/(H7+H8))"
$ws.Range("I9").Formula = "=2*((I7*I8)/(I7+I8))"
$ws.Range("J9").Formula = "=2*((J7*J8)/(J7+J8))"
$ws.Range("K9").Formula = "=2*((K7*K8)/(K7+K8))"
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 1

# --- Format the newly-filled block as percentages (matches the workbook's
#     existing K=1/K=2/K=3 block styling) -----------------------------------
$ws.Range("F7:N9").NumberFormat = "0%"

# --- Saved cursor position moved while reviewing the new numbers ----------
$ws.Range("T12").Select()
